$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59 (shifts existing rows 59-63 down to 60-64),
# mirroring Excel's "insert row above" behaviour (copies formatting
# down from the row above, e.g. the date style on column D).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Range("A59").Value2 = 8
$ws.Range("B59").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C59").Value2 = "Coquimbo"
$ws.Range("D59").Value2 = 44826
$ws.Range("E59").Value2 = 4
$ws.Range("F59").Value2 = 100114007
$ws.Range("G59").Value2 = "Jengibre"
$ws.Range("H59").Value2 = "Sin especificar"
$ws.Range("I59").Value2 = "Primera"
$ws.Range("J59").Value2 = 520
$ws.Range("K59").Value2 = 14000
$ws.Range("L59").Value2 = 15000
$ws.Range("M59").Value2 = 14500
$ws.Range("N59").Value2 = "$/caja 13 kilos"
$ws.Range("O59").Value2 = "Perú"
$ws.Range("P59").Value2 = 1115
$ws.Range("Q59").Value2 = 13
$ws.Range("R59").Value2 = "Hortaliza"
